$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.183.28"
$ws.Cells.Item(2, 5).Value = "  +0.48%  "
$ws.Cells.Item(3, 4).Value = "1.684.99"
$ws.Cells.Item(3, 5).Value = "  -0.12%  "
$ws.Cells.Item(4, 5).Value = "  +0.15%  "
$ws.Cells.Item(5, 4).Value = "'215.75"
$ws.Cells.Item(5, 5).Value = "  -0.24%  "
$ws.Cells.Item(6, 5).Value = "  +0.42%  "
$ws.Cells.Item(7, 5).Value = "  +0.20%  "
$ws.Cells.Item(8, 4).Value = "'23.10"
$ws.Cells.Item(8, 5).Value = "  +7.25%  "
$ws.Cells.Item(9, 5).Value = "  +2.78%  "
$ws.Cells.Item(10, 5).Value = "  +0.67%  "
$ws.Cells.Item(11, 5).Value = "  +0.34%  "
$ws.Cells.Item(12, 4).Value = "1.924.77"
$ws.Cells.Item(12, 5).Value = "  +0.01%  "
$ws.Cells.Item(13, 4).Value = "1.686.24"
$ws.Cells.Item(13, 5).Value = "  -0.20%  "
$ws.Cells.Item(14, 4).Value = "'4.19"
$ws.Cells.Item(14, 5).Value = "  +2.07%  "
$ws.Cells.Item(15, 5).Value = "  +3.35%  "
$ws.Cells.Item(16, 4).Value = "'66.92"
$ws.Cells.Item(16, 5).Value = "  +0.78%  "
$ws.Cells.Item(17, 4).Value = "27.185.21"
$ws.Cells.Item(17, 5).Value = "  +0.30%  "
$ws.Cells.Item(18, 4).Value = "'236.05"
$ws.Cells.Item(18, 5).Value = "  -0.83%  "
$ws.Cells.Item(19, 4).Value = "'8.01"
$ws.Cells.Item(19, 5).Value = "  -2.58%  "
$ws.Cells.Item(20, 5).Value = "  +1.15%  "
$ws.Cells.Item(21, 5).Value = "  +0.12%  "
$ws.Cells.Item(22, 5).Value = "  +2.14%  "
$ws.Cells.Item(23, 5).Value = "  +3.56%  "
$ws.Cells.Item(24, 5).Value = "  -2.70%  "
$ws.Cells.Item(25, 4).Value = "'147.41"
$ws.Cells.Item(25, 5).Value = "  +0.36%  "
$ws.Cells.Item(26, 5).Value = "  +1.10%  "
$ws.Cells.Item(27, 4).Value = "'16.45"
$ws.Cells.Item(27, 5).Value = "  +2.16%  "
$ws.Cells.Item(28, 5).Value = "  +0.51%  "
$ws.Cells.Item(29, 5).Value = "  +0.16%  "
$ws.Cells.Item(30, 4).Value = "'0.0505"
$ws.Cells.Item(30, 5).Value = "  +1.05%  "
$ws.Cells.Item(31, 5).Value = "  +0.13%  "
$ws.Cells.Item(32, 5).Value = "  +0.85%  "
$ws.Cells.Item(33, 4).Value = "1.540.52"
$ws.Cells.Item(33, 5).Value = "  +1.86%  "
$ws.Cells.Item(34, 5).Value = "  +1.58%  "
$ws.Cells.Item(35, 5).Value = "  -1.41%  "
$ws.Cells.Item(36, 5).Value = "  +2.42%  "
$ws.Cells.Item(37, 4).Value = "'0.945"
$ws.Cells.Item(37, 5).Value = "  +2.58%  "
$ws.Cells.Item(38, 5).Value = "  -0.41%  "
$ws.Cells.Item(39, 5).Value = "  -0.36%  "
$ws.Cells.Item(40, 5).Value = "  +1.49%  "
$ws.Cells.Item(41, 5).Value = "  -0.02%  "
$ws.Cells.Item(42, 4).Value = "'69.14"
$ws.Cells.Item(42, 5).Value = "  +0.80%  "
$ws.Cells.Item(43, 5).Value = "  +0.13%  "
$ws.Cells.Item(44, 5).Value = "  -1.64%  "
$ws.Cells.Item(45, 4).Value = "1.831.78"
$ws.Cells.Item(45, 5).Value = "  +0.36%  "
$ws.Cells.Item(46, 4).Value = "'0.789"
$ws.Cells.Item(46, 5).Value = "  +0.82%  "
$ws.Cells.Item(47, 4).Value = "'90.08"
$ws.Cells.Item(47, 5).Value = "  -0.30%  "
$ws.Cells.Item(48, 5).Value = "  +5.12%  "
$ws.Cells.Item(49, 5).Value = "  +5.15%  "
$ws.Cells.Item(50, 4).Value = "'8.24"
$ws.Cells.Item(50, 5).Value = "  +4.46%  "
$ws.Cells.Item(51, 5).Value = "  -1.14%  "